# Replace the coordinate table with the updated set of latitude/longitude
# points (rows 2-14), resize columns A/B to fit the new data, and move the
# active selection to A15 (the first empty row below the new table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / extend the data rows -----------------------------------
$data = @(
    @{ Row = 2;  A = 47.601369400000003;  B = -122.3232079 },
    @{ Row = 3;  A = 47.6013728;          B = -122.32346510000001 },
    @{ Row = 4;  A = 47.601076200000001;  B = -122.3227959 },
    @{ Row = 5;  A = 47.6009308;          B = -122.32296119999999 },
    @{ Row = 6;  A = 47.600820499999998;  B = -122.3225431 },
    @{ Row = 7;  A = 47.600453999999999;  B = -122.3219902 },
    @{ Row = 8;  A = 47.600075099999998;  B = -122.3218555 },
    @{ Row = 9;  A = 47.599954680000003;  B = -122.3212278 },
    @{ Row = 10; A = 47.5993377;          B = -122.3211692 },
    @{ Row = 11; A = 47.599806800000003;  B = -122.32167579999999 },
    @{ Row = 12; A = 47.600333999999997;  B = -122.32222367 },
    @{ Row = 13; A = 47.600707499999999;  B = -122.3227232 },
    @{ Row = 14; A = 47.601587799999997;  B = -122.3224958 }
)

foreach ($point in $data) {
    $ws.Cells.Item($point.Row, 1).Value = $point.A
    $ws.Cells.Item($point.Row, 2).Value = $point.B
}

# --- Column widths (A wider, B auto "best fit" sized) -----------------
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 11.83

# --- Move the selection below the new table ---------------------------
$ws.Range("A15").Select()
